$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two added date columns (V = 7-jul, W = 10-jul)
$ws.Range("V1").Value = "7-jul"
$ws.Range("W1").Value = "10-jul"

# New data values in column W (column V stays empty in the data rows,
# matching the source edit which only populated W2:W11)
$ws.Range("W2").Value = 13
$ws.Range("W3").Value = 18
$ws.Range("W4").Value = 7
$ws.Range("W5").Value = 10
$ws.Range("W6").Value = 13
$ws.Range("W7").Value = 15
$ws.Range("W8").Value = 15
$ws.Range("W9").Value = 10
$ws.Range("W10").Value = 22
$ws.Range("W11").Value = 22

# Match the formatting applied to the rest of the numeric data columns
# (integer number format, centered horizontal alignment)
$ws.Range("W2:W11").HorizontalAlignment = -4108
$ws.Range("W2:W11").NumberFormat = "0"

# Update the active selection to mirror the saved view state
$ws.Range("W5").Select() | Out-Null
